$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.307710409164429
$ws.Range("B1").Value = 3.698065042495728
$ws.Range("C1").Value = 4.104722023010254
$ws.Range("D1").Value = 2.732076406478882
$ws.Range("E1").Value = 1.051737427711487
